$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) DNA / HIV paragraph: re-save merges the "were" run + proofErr marks
#    back into the surrounding text as a single run.
# ---------------------------------------------------------------------
$find = $d.Content.Find
$dnaOld = "Developed applications to determine mutation patterns in the DNA sequences of HIV patients to assist with targeted retroviral drug therapies. Statistical analysis of large genome sequences were employed to calculate highly accurate expected ranges of mutations."
$find.Execute($dnaOld, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($find.Found) {
    $rng = $find.Parent
    $start = $rng.Start
    $len = $dnaOld.Length
    # append then remove a sentinel character: forces the engine to
    # normalise/merge the run (dropping the proofErr markers) while
    # leaving the visible text unchanged.
    $rng.Text = $dnaOld + "X"
    $d.Range($start + $len, $start + $len + 1).Delete()
}

# ---------------------------------------------------------------------
# 2) "Financial Sector" -> "financial sector" (capitalization fix)
# ---------------------------------------------------------------------
$find2 = $d.Content.Find
$find2.Execute("Financial Sector", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($find2.Found) {
    $fsRng = $find2.Parent
    $fPos = $fsRng.Start          # index of "F" in "Financial"
    $sPos = $fsRng.Start + 10     # index of "S" in "Sector"

    $d.Range($fPos, $fPos + 1).Text = "f"
    $d.Range($sPos, $sPos + 1).Text = "s"
}

Write-Output "done"
